$wb = $excel.ActiveWorkbook
Write-Host "test"
